# Auto-generated PowerShell COM-interop script applying the cryptos list update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) and Volume(1h) (column E) updates ---
$ws.Range("D2").Value = "98.759.64"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").Value = "3.336.84"
$ws.Range("E3").Value = "  -0.68%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'256.46"
$ws.Range("E5").Value = "  -2.24%  "
$ws.Range("D6").Value = "'642.31"
$ws.Range("E6").Value = "  +1.20%  "
$ws.Range("D7").Value = "'1.57"
$ws.Range("E7").Value = "  +13.13%  "
$ws.Range("D8").Value = "'0.428"
$ws.Range("E8").Value = "  +8.87%  "
$ws.Range("E9").Value = "  +25.67%  "
$ws.Range("E10").Value = "  +0.00%  "
$ws.Range("D11").Value = "3.335.56"
$ws.Range("E11").Value = "  -0.58%  "
$ws.Range("D12").Value = "'0.207"
$ws.Range("E12").Value = "  +3.18%  "
$ws.Range("D13").Value = "'43.67"
$ws.Range("E13").Value = "  +20.34%  "
$ws.Range("D14").Value = "'0.0000271"
$ws.Range("E14").Value = "  +8.93%  "
$ws.Range("D15").Value = "98.426.12"
$ws.Range("E15").Value = "  +0.04%  "
$ws.Range("D16").Value = "3.962.62"
$ws.Range("E16").Value = "  -0.03%  "
$ws.Range("E17").Value = "  -0.15%  "
$ws.Range("D18").Value = "3.336.22"
$ws.Range("E18").Value = "  +0.07%  "
$ws.Range("D19").Value = "'7.12"
$ws.Range("E19").Value = "  +15.82%  "
$ws.Range("D20").Value = "'16.62"
$ws.Range("E20").Value = "  +9.86%  "
$ws.Range("D21").Value = "'539.32"
$ws.Range("E21").Value = "  +9.05%  "
$ws.Range("D22").Value = "'3.52"
$ws.Range("E22").Value = "  -1.56%  "
$ws.Range("D23").Value = "'10.14"
$ws.Range("E23").Value = "  +8.84%  "
$ws.Range("D24").Value = "'0.443"
$ws.Range("E24").Value = "  +55.84%  "
$ws.Range("D25").Value = "'0.0000204"
$ws.Range("E25").Value = "  -3.92%  "
$ws.Range("D26").Value = "'100.98"
$ws.Range("E26").Value = "  +13.70%  "
$ws.Range("D27").Value = "'6.15"
$ws.Range("E27").Value = "  +7.78%  "
$ws.Range("D28").Value = "'12.51"
$ws.Range("E28").Value = "  +3.50%  "
$ws.Range("D29").Value = "3.513.33"
$ws.Range("E29").Value = "  -0.70%  "
$ws.Range("E30").Value = "  +17.49%  "
$ws.Range("D32").Value = "'11.03"
$ws.Range("E32").Value = "  +14.83%  "
$ws.Range("D33").Value = "'0.189"
$ws.Range("E33").Value = "  -3.87%  "
$ws.Range("E34").Value = "  +0.17%  "
$ws.Range("E35").Value = "  +4.88%  "
$ws.Range("D36").Value = "'0.525"
$ws.Range("E36").Value = "  +12.12%  "
$ws.Range("E37").Value = "  +2.64%  "
$ws.Range("D38").Value = "'2.06"
$ws.Range("E38").Value = "  +3.36%  "
$ws.Range("E39").Value = "  +2.93%  "
$ws.Range("D40").Value = "'520.90"
$ws.Range("E40").Value = "  +2.39%  "
$ws.Range("D41").Value = "'24.73"
$ws.Range("E41").Value = "  -0.45%  "
$ws.Range("E42").Value = "  +1.89%  "
$ws.Range("E43").Value = "  +2.54%  "
$ws.Range("D44").Value = "'0.811"
$ws.Range("E44").Value = "  +5.48%  "
$ws.Range("D45").Value = "'3.21"
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("D47").Value = "'0.0394"
$ws.Range("E47").Value = "  +22.51%  "
$ws.Range("E48").Value = "  +4.57%  "
$ws.Range("D51").Value = "'49.98"

# --- Row 49/50: rank swap between Cosmos and Monero ---
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").Value = "'7.75"
$ws.Range("E49").Value = "  +18.37%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").Value = "'163.71"
$ws.Range("E50").Value = "  +1.21%  "
